$wb = $excel.ActiveWorkbook

$dateFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Sheet "Overview" - add row 4 for the newly handed-back file
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A4").Value = "0fc0d813-078d-4b8a-9b21-995f72b5ad67.md"
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("E4").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F4").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G4").Value = "2016-08-29 22:47:10"
$wsOverview.Range("G4").NumberFormat = $dateFormat

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7a5195981143d3124097375ae463f5baea7f2ea/e2e/0fc0d813-078d-4b8a-9b21-995f72b5ad67.md",
    "",
    "",
    "e2e\0fc0d813-078d-4b8a-9b21-995f72b5ad67.md"
) | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn" - add row 4 for the newly handed-back file
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Range("B4").Value = ".md"
$wsZhCn.Range("C4").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("D4").Value = "e2e"
$wsZhCn.Range("E4").Value = "ht"
$wsZhCn.Range("F4").Value = "'True"
$wsZhCn.Range("G4").Value = "0fc0d813-078d-4b8a-9b21-995f72b5ad67.32d64c440a37303dbc5c1a203e920ac94105ad1c.zh-cn.xlf"
$wsZhCn.Range("H4").Value = "2016-08-29 22:47:00"
$wsZhCn.Range("H4").NumberFormat = $dateFormat
$wsZhCn.Range("J4").Value = "0fc0d813-078d-4b8a-9b21-995f72b5ad67.32d64c440a37303dbc5c1a203e920ac94105ad1c.zh-cn.xlf"
$wsZhCn.Range("K4").Value = "2016-08-29 22:47:29"
$wsZhCn.Range("K4").NumberFormat = $dateFormat
$wsZhCn.Range("L4").Value = "'"
$wsZhCn.Range("M4").Value = "'True"
$wsZhCn.Range("N4").Value = "'"
$wsZhCn.Range("O4").Value = "'False"
$wsZhCn.Range("P4").Value = "'"

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7a5195981143d3124097375ae463f5baea7f2ea/e2e/0fc0d813-078d-4b8a-9b21-995f72b5ad67.md",
    "",
    "",
    "0fc0d813-078d-4b8a-9b21-995f72b5ad67.md"
) | Out-Null

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("I4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/e22c739602a20992365fb53ff511b8b6310c9095/e2e/0fc0d813-078d-4b8a-9b21-995f72b5ad67.md",
    "",
    "",
    "0fc0d813-078d-4b8a-9b21-995f72b5ad67.md"
) | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de" - add row 4 for the newly handed-back file
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Range("B4").Value = ".md"
$wsDeDe.Range("C4").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("D4").Value = "e2e"
$wsDeDe.Range("E4").Value = "ht"
$wsDeDe.Range("F4").Value = "'True"
$wsDeDe.Range("G4").Value = "0fc0d813-078d-4b8a-9b21-995f72b5ad67.32d64c440a37303dbc5c1a203e920ac94105ad1c.de-de.xlf"
$wsDeDe.Range("H4").Value = "2016-08-29 22:47:10"
$wsDeDe.Range("H4").NumberFormat = $dateFormat
$wsDeDe.Range("J4").Value = "0fc0d813-078d-4b8a-9b21-995f72b5ad67.32d64c440a37303dbc5c1a203e920ac94105ad1c.de-de.xlf"
$wsDeDe.Range("K4").Value = "2016-08-29 22:47:36"
$wsDeDe.Range("K4").NumberFormat = $dateFormat
$wsDeDe.Range("L4").Value = "'"
$wsDeDe.Range("M4").Value = "'True"
$wsDeDe.Range("N4").Value = "'"
$wsDeDe.Range("O4").Value = "'False"
$wsDeDe.Range("P4").Value = "'"

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7a5195981143d3124097375ae463f5baea7f2ea/e2e/0fc0d813-078d-4b8a-9b21-995f72b5ad67.md",
    "",
    "",
    "0fc0d813-078d-4b8a-9b21-995f72b5ad67.md"
) | Out-Null

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("I4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/30f7d4044120650e46846d7f69110194ba2eb1a/e2e/0fc0d813-078d-4b8a-9b21-995f72b5ad67.md",
    "",
    "",
    "0fc0d813-078d-4b8a-9b21-995f72b5ad67.md"
) | Out-Null
